$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as numbers by Excel
# are pre-formatted as Text so they remain exact strings (matches the
# inlineStr text cells in the source workbook).
$textCells = @("D5","D8","D9","D10","D11","D18","D19","D21","D24","D25","D26","D27","D33","D40","D42","D43","D44","D45","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.563.80"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.810.23"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "226.08"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  +3.29%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "37.55"
$ws.Range("E8").Value = "  +7.41%  "
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  -3.10%  "
$ws.Range("D10").Value = "0.0681"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "2.071.62"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "1.808.58"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("D16").Value = "34.536.54"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "68.60"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "244.11"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "11.23"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("D25").Value = "172.10"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "7.87"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "17.41"
$ws.Range("E27").Value = "  +3.44%  "
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").Value = "0.0521"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "1.366.36"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("E38").Value = "  -4.15%  "
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "1.22"
$ws.Range("E40").Value = "  +8.37%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.78"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "81.01"
$ws.Range("E43").Value = "  -3.48%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.941"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").Value = "14.01"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").Value = "1.972.54"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "103.04"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "0.0₆0121"
$ws.Range("E51").Value = "  -7.73%  "
